$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: new match (LATVIA - VIRSLIGA, Metta vs RFS) ---
$ws.Range("A3").Value = "pbRMNv36"
$ws.Range("B3").Value = "28/10/2024"
$ws.Range("C3").Value = "13:30"
$ws.Range("D3").Value = "LATVIA - VIRSLIGA"
$ws.Range("E3").Value = "Metta"
$ws.Range("F3").Value = "RFS"
$ws.Range("G3").Value = 10.75
$ws.Range("H3").Value = 5.7
$ws.Range("I3").Value = 1.19
$ws.Range("J3").Value = 8.25
$ws.Range("K3").Value = 2.77
$ws.Range("L3").Value = 1.55
$ws.Range("M3").Value = 1.01
$ws.Range("N3").Value = 18.5
$ws.Range("O3").Value = 1.06
$ws.Range("P3").Value = 5.6
$ws.Range("Q3").Value = 1.33
$ws.Range("R3").Value = 2.73
$ws.Range("S3").Value = 1.19
$ws.Range("T3").Value = 4.15
$ws.Range("U3").Value = 1.87
$ws.Range("V3").Value = 1.89
$ws.Range("W3").Value = 30
$ws.Range("X3").Value = 75
$ws.Range("Y3").Value = 28
$ws.Range("Z3").Value = 250
$ws.Range("AA3").Value = 100
$ws.Range("AB3").Value = 65
$ws.Range("AC3").Value = 19
$ws.Range("AD3").Value = 11
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 65
$ws.Range("AG3").Value = 350
$ws.Range("AH3").Value = 8.25
$ws.Range("AI3").Value = 6.2
$ws.Range("AJ3").Value = 8
$ws.Range("AK3").Value = 6.5
$ws.Range("AL3").Value = 8.5
$ws.Range("AM3").Value = 20
$ws.Range("AN3").Value = 11.5
$ws.Range("AO3").Value = 60
$ws.Range("AP3").Value = 45
$ws.Range("AQ3").Value = 450
$ws.Range("AR3").Value = 350
$ws.Range("AS3").Value = 500
$ws.Range("AT3").Value = 3.9
$ws.Range("AU3").Value = 8.75
$ws.Range("AV3").Value = 65
$ws.Range("AW3").Value = 3.25
$ws.Range("AX3").Value = 4.9
$ws.Range("AY3").Value = 13.5
$ws.Range("AZ3").Value = 10.75
$ws.Range("BA3").Value = 30
$ws.Range("BB3").Value = 150
$ws.Range("BC3").Value = 51
$ws.Range("BD3").Value = 51

# --- Row 4: new match (ROMANIA - LIGA 1, FC Botosani vs U. Cluj) ---
$ws.Range("A4").Value = "b7I8pNte"
$ws.Range("B4").Value = "28/10/2024"
$ws.Range("C4").Value = "13:00"
$ws.Range("D4").Value = "ROMANIA - LIGA 1"
$ws.Range("E4").Value = "FC Botosani"
$ws.Range("F4").Value = "U. Cluj"
$ws.Range("G4").Value = 3.6
$ws.Range("H4").Value = 3
$ws.Range("I4").Value = 2.15
$ws.Range("J4").Value = 4.33
$ws.Range("K4").Value = 1.95
$ws.Range("L4").Value = 3
$ws.Range("M4").Value = 1.1
$ws.Range("N4").Value = 7
$ws.Range("O4").Value = 1.5
$ws.Range("P4").Value = 2.5
$ws.Range("Q4").Value = 2.5
$ws.Range("R4").Value = 1.5
$ws.Range("S4").Value = 1.53
$ws.Range("T4").Value = 2.38
$ws.Range("U4").Value = 2.1
$ws.Range("V4").Value = 1.67
$ws.Range("W4").Value = 8.5
$ws.Range("X4").Value = 17
$ws.Range("Y4").Value = 13
$ws.Range("Z4").Value = 41
$ws.Range("AA4").Value = 34
$ws.Range("AB4").Value = 41
$ws.Range("AC4").Value = 6.5
$ws.Range("AD4").Value = 6
$ws.Range("AE4").Value = 19
$ws.Range("AF4").Value = 67
$ws.Range("AG4").Value = 1000
$ws.Range("AH4").Value = 6
$ws.Range("AI4").Value = 9
$ws.Range("AJ4").Value = 10
$ws.Range("AK4").Value = 19
$ws.Range("AL4").Value = 21
$ws.Range("AM4").Value = 41
$ws.Range("AN4").Value = 5.5
$ws.Range("AO4").Value = 21
$ws.Range("AP4").Value = 34
$ws.Range("AQ4").Value = 81
$ws.Range("AR4").Value = 126
$ws.Range("AS4").Value = 351
$ws.Range("AT4").Value = 2.38
$ws.Range("AU4").Value = 9
$ws.Range("AV4").Value = 67
$ws.Range("AW4").Value = 4
$ws.Range("AX4").Value = 13
$ws.Range("AY4").Value = 29
$ws.Range("AZ4").Value = 41
$ws.Range("BA4").Value = 81
$ws.Range("BB4").Value = 251
$ws.Range("BC4").Value = 51
$ws.Range("BD4").Value = 51

# --- Row 5: update some odds values (TURKEY - SUPER LIG stays) ---
$ws.Range("G5").Value = 2.1
$ws.Range("I5").Value = 3.4
$ws.Range("L5").Value = 3.75
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 3
$ws.Range("U5").Value = 1.67
$ws.Range("V5").Value = 2.1
$ws.Range("W5").Value = 8.5
$ws.Range("X5").Value = 11
$ws.Range("AC5").Value = 11
$ws.Range("AG5").Value = 151
$ws.Range("AH5").Value = 12
$ws.Range("AJ5").Value = 12
$ws.Range("AL5").Value = 26
$ws.Range("AN5").Value = 4.33
$ws.Range("AS5").Value = 126
$ws.Range("AT5").Value = 3
$ws.Range("AU5").Value = 7.5
